# Small fix: standardised cost_variable_om
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Parameter column (C) for rows 10-39 from "cost_variable" to "cost_variable_om"
$range = $ws.Range("C10:C39")
$range.Value = "cost_variable_om"

# Mirror the selection change recorded in the saved view state
$range.Select()
